$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.2855834829871 ; $ws.Range("C2").Value = 11.21434405186043 ; $ws.Range("D2").Value = 9.880341846853312 ; $ws.Range("E2").Value = 13.9544573197214 ; $ws.Range("F2").Value = 29.41981128083354 ; $ws.Range("H2").Value = 7.344005520526261 ; $ws.Range("I2").Value = 19.27409676538315 ; $ws.Range("J2").Value = 9.787616231031418 ; $ws.Range("O2").Value = 21.72566271065677
$ws.Range("B3").Value = 15.53553635930213 ; $ws.Range("C3").Value = 10.57202405173073 ; $ws.Range("D3").Value = 9.811964237362123 ; $ws.Range("E3").Value = 13.88622203461771 ; $ws.Range("F3").Value = 29.50437966604834 ; $ws.Range("H3").Value = 7.344005520526261 ; $ws.Range("I3").Value = 19.44544583590744 ; $ws.Range("J3").Value = 9.795528358734048 ; $ws.Range("O3").Value = 21.84320915618729
$ws.Range("B4").Value = 15.05590403970333 ; $ws.Range("C4").Value = 10.15638819930145 ; $ws.Range("D4").Value = 9.771021521452878 ; $ws.Range("E4").Value = 13.84691397622056 ; $ws.Range("F4").Value = 29.56637644760508 ; $ws.Range("H4").Value = 7.344005520526261 ; $ws.Range("I4").Value = 19.5562783632397 ; $ws.Range("J4").Value = 9.802050078544365 ; $ws.Range("O4").Value = 21.92237655212764
$ws.Range("B5").Value = 14.85587090174414 ; $ws.Range("C5").Value = 9.981770028799989 ; $ws.Range("D5").Value = 9.75461165354373 ; $ws.Range("E5").Value = 13.83155875148785 ; $ws.Range("F5").Value = 29.59416026829772 ; $ws.Range("H5").Value = 7.344005520526261 ; $ws.Range("I5").Value = 19.60285896525389 ; $ws.Range("J5").Value = 9.805126129648992 ; $ws.Range("O5").Value = 21.95638776736269
$ws.Range("B6").Value = 14.82238620963225 ; $ws.Range("C6").Value = 9.952461515921398 ; $ws.Range("D6").Value = 9.751903758068543 ; $ws.Range("E6").Value = 13.82904941785821 ; $ws.Range("F6").Value = 29.5989255335807 ; $ws.Range("H6").Value = 7.344005520526261 ; $ws.Range("I6").Value = 19.61067917040405 ; $ws.Range("J6").Value = 9.805662175315506 ; $ws.Range("O6").Value = 21.96214073579839
$ws.Range("B7").Value = 15.05322454225488 ; $ws.Range("C7").Value = 10.15405431818319 ; $ws.Range("D7").Value = 9.770799083773646 ; $ws.Range("E7").Value = 13.84670418987355 ; $ws.Range("F7").Value = 29.56674096587966 ; $ws.Range("H7").Value = 7.344005520526261 ; $ws.Range("I7").Value = 19.55690083258499 ; $ws.Range("J7").Value = 9.802089869175624 ; $ws.Range("O7").Value = 21.92282816572744
$ws.Range("B8").Value = 16.03105545504948 ; $ws.Range("C8").Value = 10.99733553633772 ; $ws.Range("D8").Value = 9.856556874288749 ; $ws.Range("E8").Value = 13.93040008656549 ; $ws.Range("F8").Value = 29.44687265140232 ; $ws.Range("H8").Value = 7.344005520526261 ; $ws.Range("I8").Value = 19.3320115771992 ; $ws.Range("J8").Value = 9.789999111785546 ; $ws.Range("O8").Value = 21.76473603904757
$ws.Range("B9").Value = 17.78882160445051 ; $ws.Range("C9").Value = 12.47904981584823 ; $ws.Range("D9").Value = 10.03241323302142 ; $ws.Range("E9").Value = 14.11449281049924 ; $ws.Range("F9").Value = 29.29225653723545 ; $ws.Range("H9").Value = 7.344005520526261 ; $ws.Range("I9").Value = 18.93556228612061 ; $ws.Range("J9").Value = 9.779483241160852 ; $ws.Range("O9").Value = 21.51059287121511
$ws.Range("B10").Value = 18.97385515298169 ; $ws.Range("C10").Value = 13.45949919669951 ; $ws.Range("D10").Value = 10.16550534432579 ; $ws.Range("E10").Value = 14.26109401486267 ; $ws.Range("F10").Value = 29.22836151625194 ; $ws.Range("H10").Value = 7.344005520526261 ; $ws.Range("I10").Value = 18.67136900038363 ; $ws.Range("J10").Value = 9.779787126839807 ; $ws.Range("O10").Value = 21.35845922945906
$ws.Range("B11").Value = 19.48843426857424 ; $ws.Range("C11").Value = 13.88159880988164 ; $ws.Range("D11").Value = 10.22672274331498 ; $ws.Range("E11").Value = 14.33006039955095 ; $ws.Range("F11").Value = 29.21020329807997 ; $ws.Range("H11").Value = 7.344005520526261 ; $ws.Range("I11").Value = 18.55704456381937 ; $ws.Range("J11").Value = 9.781663632643593 ; $ws.Range("O11").Value = 21.29686976760731
$ws.Range("B12").Value = 19.67967655310752 ; $ws.Range("C12").Value = 14.03797336093682 ; $ws.Range("D12").Value = 10.24998601723073 ; $ws.Range("E12").Value = 14.35648653858987 ; $ws.Range("F12").Value = 29.20490337592358 ; $ws.Range("H12").Value = 7.344005520526261 ; $ws.Range("I12").Value = 18.51459456962926 ; $ws.Range("J12").Value = 9.782623474346131 ; $ws.Range("O12").Value = 21.27465169662584
$ws.Range("B13").Value = 19.63865133637584 ; $ws.Range("C13").Value = 14.00444977722254 ; $ws.Range("D13").Value = 10.24497247008266 ; $ws.Range("E13").Value = 14.35078168794683 ; $ws.Range("F13").Value = 29.20597459952263 ; $ws.Range("H13").Value = 7.344005520526261 ; $ws.Range("I13").Value = 18.52369948099734 ; $ws.Range("J13").Value = 9.782405682567907 ; $ws.Range("O13").Value = 21.27938748656658
$ws.Range("B14").Value = 19.5042408739799 ; $ws.Range("C14").Value = 13.89453341762394 ; $ws.Range("D14").Value = 10.22863507945414 ; $ws.Range("E14").Value = 14.33222838085769 ; $ws.Range("F14").Value = 29.20973564183904 ; $ws.Range("H14").Value = 7.344005520526261 ; $ws.Range("I14").Value = 18.55353529936007 ; $ws.Range("J14").Value = 9.781737608613616 ; $ws.Range("O14").Value = 21.29501968497905
$ws.Range("B15").Value = 19.42143692880061 ; $ws.Range("C15").Value = 13.82675450582887 ; $ws.Range("D15").Value = 10.21863812618598 ; $ws.Range("E15").Value = 14.32090380094955 ; $ws.Range("F15").Value = 29.21224486089154 ; $ws.Range("H15").Value = 7.344005520526261 ; $ws.Range("I15").Value = 18.57192027292681 ; $ws.Range("J15").Value = 9.781360830296999 ; $ws.Range("O15").Value = 21.30473897039227
$ws.Range("B16").Value = 18.93972473454828 ; $ws.Range("C16").Value = 13.43143034801721 ; $ws.Range("D16").Value = 10.16151688917439 ; $ws.Range("E16").Value = 14.256631244649 ; $ws.Range("F16").Value = 29.22976837890435 ; $ws.Range("H16").Value = 7.344005520526261 ; $ws.Range("I16").Value = 18.67895823295019 ; $ws.Range("J16").Value = 9.779699413437498 ; $ws.Range("O16").Value = 21.36263830853166
$ws.Range("B17").Value = 18.63786233961212 ; $ws.Range("C17").Value = 13.18276571476331 ; $ws.Range("D17").Value = 10.12663646429651 ; $ws.Range("E17").Value = 14.21777309098722 ; $ws.Range("F17").Value = 29.24331814355063 ; $ws.Range("H17").Value = 7.344005520526261 ; $ws.Range("I17").Value = 18.74612281110439 ; $ws.Range("J17").Value = 9.779124936609231 ; $ws.Range("O17").Value = 21.4001151519464
$ws.Range("B18").Value = 18.46193765244654 ; $ws.Range("C18").Value = 13.03749076478738 ; $ws.Range("D18").Value = 10.10663844297293 ; $ws.Range("E18").Value = 14.19563828240581 ; $ws.Range("F18").Value = 29.25213794290246 ; $ws.Range("H18").Value = 7.344005520526261 ; $ws.Range("I18").Value = 18.78530552638712 ; $ws.Range("J18").Value = 9.778958177653395 ; $ws.Range("O18").Value = 21.42238723592448
$ws.Range("B19").Value = 18.40198047193277 ; $ws.Range("C19").Value = 12.98791773696131 ; $ws.Range("D19").Value = 10.09987895845285 ; $ws.Range("E19").Value = 14.18818132626 ; $ws.Range("F19").Value = 29.25530019418609 ; $ws.Range("H19").Value = 7.344005520526261 ; $ws.Range("I19").Value = 18.79866683862836 ; $ws.Range("J19").Value = 9.778929845008431 ; $ws.Range("O19").Value = 21.43005093916882
$ws.Range("B20").Value = 18.67023503511529 ; $ws.Range("C20").Value = 13.20946953730352 ; $ws.Range("D20").Value = 10.1303430026999 ; $ws.Range("E20").Value = 14.2218874412512 ; $ws.Range("F20").Value = 29.24176946662718 ; $ws.Range("H20").Value = 7.344005520526261 ; $ws.Range("I20").Value = 18.73891595877174 ; $ws.Range("J20").Value = 9.779169156665846 ; $ws.Range("O20").Value = 21.39605146651634
$ws.Range("B21").Value = 19.54381934076168 ; $ws.Range("C21").Value = 13.9269127441387 ; $ws.Range("D21").Value = 10.23343167489646 ; $ws.Range("E21").Value = 14.33766966249883 ; $ws.Range("F21").Value = 29.20858810059255 ; $ws.Range("H21").Value = 7.344005520526261 ; $ws.Range("I21").Value = 18.54474894066055 ; $ws.Range("J21").Value = 9.781927079896567 ; $ws.Range("O21").Value = 21.29039808020725
$ws.Range("B22").Value = 20.09363790688439 ; $ws.Range("C22").Value = 14.37559686338574 ; $ws.Range("D22").Value = 10.30127445512182 ; $ws.Range("E22").Value = 14.41513894711356 ; $ws.Range("F22").Value = 29.19609205512402 ; $ws.Range("H22").Value = 7.344005520526261 ; $ws.Range("I22").Value = 18.42275849307607 ; $ws.Range("J22").Value = 9.785181935131321 ; $ws.Range("O22").Value = 21.22778965276528
$ws.Range("B23").Value = 19.80214895848603 ; $ws.Range("C23").Value = 14.13798159701823 ; $ws.Range("D23").Value = 10.26502767497743 ; $ws.Range("E23").Value = 14.37363338836159 ; $ws.Range("F23").Value = 29.2019183758708 ; $ws.Range("H23").Value = 7.344005520526261 ; $ws.Range("I23").Value = 18.48741796266713 ; $ws.Range("J23").Value = 9.783312133247588 ; $ws.Range("O23").Value = 21.26061248202446
$ws.Range("B24").Value = 18.65560674522315 ; $ws.Range("C24").Value = 13.1974039457414 ; $ws.Range("D24").Value = 10.12866710360462 ; $ws.Range("E24").Value = 14.22002670283554 ; $ws.Range("F24").Value = 29.24246641606288 ; $ws.Range("H24").Value = 7.344005520526261 ; $ws.Range("I24").Value = 18.74217240451517 ; $ws.Range("J24").Value = 9.779148655389804 ; $ws.Range("O24").Value = 21.39788639743489
$ws.Range("B25").Value = 17.33143084942772 ; $ws.Range("C25").Value = 12.09701482776442 ; $ws.Range("D25").Value = 9.984097558982874 ; $ws.Range("E25").Value = 14.06263536441634 ; $ws.Range("F25").Value = 29.32539944408653 ; $ws.Range("H25").Value = 7.344005520526261 ; $ws.Range("I25").Value = 19.03805146126554 ; $ws.Range("J25").Value = 9.7809164476151 ; $ws.Range("O25").Value = 21.57330925135415
